$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8570482730865479
$ws.Range("B1").Value = 1.321430206298828
$ws.Range("D1").Value = 1.723708391189575
$ws.Range("E1").Value = 1.12641167640686
